$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.996.29'
$ws.Range("E2").Value = '  -0.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.867.80'
$ws.Range("E3").Value = '  -2.88%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.88'
$ws.Range("E5").Value = '  -1.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5067'
$ws.Range("E7").Value = '  -1.62%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3904'
$ws.Range("E8").Value = '  -2.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08148'
$ws.Range("E9").Value = '  -3.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.98'
$ws.Range("E10").Value = '  -2.20%  '

$ws.Range("E11").Value = '  -2.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.55'
$ws.Range("E12").Value = '  +6.70%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.871.39'
$ws.Range("E13").Value = '  -2.57%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.255'
$ws.Range("E14").Value = '  -0.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.141'
$ws.Range("E15").Value = '  -2.73%  '

$ws.Range("E16").Value = '  +0.05%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.53'
$ws.Range("E17").Value = '  -2.99%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001072'
$ws.Range("E18").Value = '  -3.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06332'
$ws.Range("E19").Value = '  -6.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.77'
$ws.Range("E20").Value = '  -1.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.985.11'
$ws.Range("E22").Value = '  -0.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.775'
$ws.Range("E23").Value = '  -4.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.03'
$ws.Range("E24").Value = '  -1.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.204'
$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.089.86'
$ws.Range("E26").Value = '  -2.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.73'
$ws.Range("E27").Value = '  +0.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.77'
$ws.Range("E28").Value = '  -0.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.207'
$ws.Range("E29").Value = '  -10.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.21'
$ws.Range("E30").Value = '  -2.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1030'
$ws.Range("E31").Value = '  -2.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.036'
$ws.Range("E32").Value = '  -3.52%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.854'
$ws.Range("E33").Value = '  -3.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.721'
$ws.Range("E34").Value = '  +1.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02408'
$ws.Range("E35").Value = '  -3.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.168'
$ws.Range("E36").Value = '  -0.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06309'
$ws.Range("E37").Value = '  -4.12%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2131'
$ws.Range("E38").Value = '  -3.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.167'
$ws.Range("E39").Value = '  -5.88%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.476'
$ws.Range("E40").Value = '  -5.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.208'
$ws.Range("E41").Value = '  -2.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6246'
$ws.Range("E42").Value = '  -4.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.16'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  -0.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5857'
$ws.Range("E45").Value = '  -4.31%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.80'
$ws.Range("E46").Value = '  -2.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.625'
$ws.Range("E47").Value = '  -2.69%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.981'
$ws.Range("E48").Value = '  -3.38%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.84'
$ws.Range("E49").Value = '  -2.57%  '

$ws.Range("E50").Value = '  -3.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.119'
$ws.Range("E51").Value = '  -2.95%  '
